# ---------------------------------------------------------------------------
# "code moved to khushboo folder"
#
# Semantic changes applied (derived from the canonical-OOXML diff):
#   1. Para "Copy all files from project folder directly to htdocs folder"
#        -> "Copy entire  khushboo folder directly to htdocs folder"
#   2. The two list items
#        "Change the path of project and db configuration in /config/config.php file"
#        "Change the path in  /js/script.js file"
#      are removed entirely (that step is now obsolete once the project lives
#      directly in the khushboo folder).
#   3. The API url "http://localhost/restApi/" becomes
#      "http://localhost/khushboo/restApi/".
#   4. Word's auto-managed "_GoBack" (last-edit) bookmark follows the final
#      edit location, i.e. moves from the "dbSql" list item to the end of the
#      (now edited) first list item.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Detach the existing "_GoBack" bookmark -----------------------------
# It currently sits inside the "Import db queries from dbSql folder" item;
# it will be re-created at the new last-edit location below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Remove the two obsolete "Change the path ..." list items -----------
$pChangePath1 = $d.Paragraphs.Item(3)
$pChangePath2 = $d.Paragraphs.Item(4)
$obsoleteRange = $d.Range($pChangePath1.Range.Start, $pChangePath2.Range.End)
$obsoleteRange.Delete()

# --- 3. Rewrite the "Copy ..." list item ------------------------------------
# "Copy all files from project folder ..." -> "Copy entire  khushboo folder ..."
$d.Content.Find.Execute("all files from ", $false, $false, $false, $false, `
    $false, $true, 1, $false, "entire  ", 2) | Out-Null
$d.Content.Find.Execute("project", $true, $false, $false, $false, `
    $false, $true, 1, $false, "khushboo", 2) | Out-Null

# --- 4. Re-create "_GoBack" at the end of that item (after "folder") -------
$pCopy = $d.Paragraphs.Item(2)
$endPos = $pCopy.Range.End - 1   # just before the paragraph mark

# Placing a zero-length bookmark directly at a paragraph-end boundary is
# unreliable, so nudge it into place: insert a throwaway character, wrap the
# bookmark around it, then delete the character -- the bookmark collapses
# back to the exact insertion point (standard bookmark "gravity").
$d.Range($endPos, $endPos).InsertBefore("x")
$markerRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)
$d.Range($endPos, $endPos + 1).Delete()

# --- 5. Update the REST API url --------------------------------------------
$d.Content.Find.Execute("http://localhost/restApi/", $true, $false, $false, `
    $false, $false, $true, 1, $false, "http://localhost/khushboo/restApi/", 2) | Out-Null

Write-Output "done"
